$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B7").Interior.Color = 13561798
